$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 1003
$ws.Range("B5").Value = 0.64304398148148145
$ws.Range("B5").NumberFormat = "h:mm:ss"
$ws.Range("C5").Value = 0.25

$ws.Range("C5").Select()
